$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("СВОДНАЯ")

$ws.Range("C18").Value = "Очень много простой, но рутинной работы"

$ws.Range("C22").Value = "Сущенствует много способов решения одной задачи"
$ws.Range("C22").WrapText = $true
$ws.Rows(22).RowHeight = 30

$ws.Range("C26").Value = "Хэш таблицы быстрее:)"

$ws.Range("C27").Select()
